$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new product ("Tapas para empanadas") was inserted as row 3 of the sheet,
# pushing every existing data row (previously rows 3..34) one row down
# (now rows 4..35). Re-create that by copying whole rows from the bottom up
# (so we never overwrite a row before it has been copied), then fill the
# freed-up row 3 with the new record's data.
# ---------------------------------------------------------------------------

for ($r = 34; $r -ge 3; $r--) {
    $src = $ws.Range("A" + $r + ":O" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":O" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4163)  # xlPasteValues
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}
$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# Populate the new row 3 with the "Tapas para empanadas" article.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = 7798118140024
$ws.Range("B3").Value2 = "Tapas"
$ws.Range("C3").Value2 = "para empanadas"
$ws.Range("D3").Value2 = '"freir"'
$ws.Range("E3").Value2 = "Di Pascualle"
$ws.Range("F3").Value2 = 12
$ws.Range("G3").Value2 = "und."
$ws.Range("H3").Value2 = "Bolsa"
$ws.Range("I3").Value2 = "Pastas frescas"
$ws.Range("J3").Value2 = "Argentina"
$ws.Range("K3").Value2 = 30
$ws.Range("L3").Value2 = $false
$ws.Range("M3").Value2 = $true
$ws.Range("N3").Value2 = "C:\VentaSoft\Imágenes de artículos\7798118140024.png"
$ws.Range("O3").Value2 = $true

# Make sure the new row carries the same cell styles as the rest of the
# table (numeric-id format on Codigo, bordered/boolean style on
# ImagenExactaDelArticulo, matching row 2 immediately above it), again via
# a formats-only paste so the existing style entries are reused.
$ws.Range("A2:O2").Copy()
$ws.Range("A3:O3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()
